# Add "edge 0-shot" summary rows: sample standard deviation (row 106)
# and 95% confidence-interval half-width (row 107) for columns B:F,
# mirroring the existing AVERAGE row (row 105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 106: STDEV.S over the same data range used by the AVERAGE row (105).
$ws.Range("B106").Formula = "=STDEV.S(B2:B104)"
$ws.Range("C106:F106").Formula = "=STDEV.S(C2:C104)"

# Row 107: 95% CI half-width = stdev / sqrt(n-1) * 1.96, n = 104 rows of data (103 here).
$ws.Range("B107").Formula = "=B106/SQRT(103)*1.96"
$ws.Range("C107:F107").Formula = "=C106/SQRT(103)*1.96"

# Scroll the view down and land the selection on J106, like the saved workbook.
$excel.ActiveWindow.ScrollRow = 82
[void]$ws.Range("J106").Select()
